{"js": "// Update the counts / percentages in the \"cat / n / Arrested\" table to\n// reflect the refreshed 2023-06-20 quarterly pull.\nconst replacements = [\n  [\"349\", \"368\"],\n  [\"157 (45.0)\", \"164 (44.6)\"],\n  [\"106\", \"113\"],\n  [\"47 (44.3)\", \"49 (43.4)\"],\n  [\"109\", \"114\"],\n  [\"51 (46.8)\", \"53 (46.5)\"],\n  [\"107\", \"112\"],\n  [\"43 (40.2)\", \"46 (41.1)\"],\n  [\"27\", \"29\"],\n  [\"16 (59.3)\", \"16 (55.2)\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the counts / percentages in the \"cat / n / Arrested\" table to\n# reflect the refreshed 2023-06-20 quarterly pull.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"349\", \"368\"),\n  @(\"157 (45.0)\", \"164 (44.6)\"),\n  @(\"106\", \"113\"),\n  @(\"47 (44.3)\", \"49 (43.4)\"),\n  @(\"109\", \"114\"),\n  @(\"51 (46.8)\", \"53 (46.5)\"),\n  @(\"107\", \"112\"),\n  @(\"43 (40.2)\", \"46 (41.1)\"),\n  @(\"27\", \"29\"),\n  @(\"16 (59.3)\", \"16 (55.2)\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n\n  $find.Execute($oldText, $true, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
